$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1272.8572
$ws.Range("J112").Value = 1272.5
$ws.Range("L112").Value = 3817.5
$ws.Range("N112").Value = -6033.5
$ws.Range("H138").Value = 3657.0962
$ws.Range("J138").Value = 3952.3416
$ws.Range("L138").Value = 11857.0248
$ws.Range("N138").Value = -22137.0248

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1735.0714
$ws.Range("I2").Value = 1735.0714
$ws.Range("K2").Value = 1735.0714
$ws.Range("M2").Value = -1622.0714
$ws.Range("H32").Value = 1627.966
$ws.Range("I32").Value = 1004.7778
$ws.Range("K32").Value = 1004.7778
$ws.Range("M32").Value = -717.7778
$ws.Range("H43").Value = 136106
$ws.Range("I43").Value = 222222
$ws.Range("J43").Value = 49990
$ws.Range("K43").Value = 222222
$ws.Range("L43").Value = 49990
$ws.Range("M43").Value = -221909
$ws.Range("N43").Value = -50616
$ws.Range("H116").Value = 1735.0714
$ws.Range("I116").Value = 1735.0714
$ws.Range("K116").Value = 1735.0714
$ws.Range("M116").Value = 558.9286
$ws.Range("H122").Value = 1332.826
$ws.Range("I122").Value = 1241.1177
$ws.Range("J122").Value = 1592.6666
$ws.Range("K122").Value = 3723.3531
$ws.Range("L122").Value = 4777.9998
$ws.Range("M122").Value = -1273.3531
$ws.Range("N122").Value = -9677.9998
$ws.Range("H132").Value = 3352.2
$ws.Range("I132").Value = 3122.4
$ws.Range("J132").Value = 4041.6
$ws.Range("K132").Value = 9367.200000000001
$ws.Range("L132").Value = 12124.8
$ws.Range("M132").Value = -6837.200000000001
$ws.Range("N132").Value = -17184.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1735.0714
$ws.Range("I3").Value = 1735.0714
$ws.Range("K3").Value = 1735.0714
$ws.Range("M3").Value = -1621.0714
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H94").Value = 1064.0588
$ws.Range("J94").Value = 1150.8889
$ws.Range("L94").Value = 1150.8889
$ws.Range("N94").Value = -2052.8889
$ws.Range("H96").Value = 30000
$ws.Range("I96").Value = 30000
$ws.Range("K96").Value = 30000
$ws.Range("M96").Value = -27254
$ws.Range("H105").Value = 1524.4
$ws.Range("I105").Value = 1603.2858
$ws.Range("J105").Value = 1340.3334
$ws.Range("K105").Value = 1603.2858
$ws.Range("L105").Value = 1340.3334
$ws.Range("M105").Value = 143.7141999999999
$ws.Range("N105").Value = -4834.3334
$ws.Range("H134").Value = 1522.68
$ws.Range("I134").Value = 981.96106
$ws.Range("K134").Value = 2945.88318
$ws.Range("M134").Value = -410.8831799999998

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4079.6206
$ws.Range("I31").Value = 1675.85
$ws.Range("K31").Value = 1675.85
$ws.Range("M31").Value = -1380.85
$ws.Range("H34").Value = 4079.6206
$ws.Range("I34").Value = 1675.85
$ws.Range("K34").Value = 1675.85
$ws.Range("M34").Value = -1473.85
$ws.Range("H76").Value = 7153
$ws.Range("I76").Value = 7153
$ws.Range("K76").Value = 7153
$ws.Range("M76").Value = -6838
$ws.Range("H79").Value = 7153
$ws.Range("I79").Value = 7153
$ws.Range("K79").Value = 7153
$ws.Range("M79").Value = -6061
$ws.Range("H132").Value = 1736.2941
$ws.Range("I132").Value = 1736.2941
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5208.8823
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2678.8823
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 8281.754999999999
$ws.Range("I134").Value = 8495.341
$ws.Range("K134").Value = 25486.023
$ws.Range("M134").Value = -22951.023

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 3623.5
$ws.Range("J62").Value = 3998
$ws.Range("L62").Value = 11994
$ws.Range("N62").Value = -13366
$ws.Range("H65").Value = 3623.5
$ws.Range("J65").Value = 3998
$ws.Range("L65").Value = 35982
$ws.Range("N65").Value = -42846
$ws.Range("H124").Value = 2717
$ws.Range("I124").Value = 2717
$ws.Range("K124").Value = 8151
$ws.Range("M124").Value = -3241
$ws.Range("H125").Value = 15620.25
$ws.Range("I125").Value = 4110
$ws.Range("J125").Value = 22526.4
$ws.Range("K125").Value = 12330
$ws.Range("L125").Value = 67579.20000000001
$ws.Range("M125").Value = -7410
$ws.Range("N125").Value = -77419.20000000001
$ws.Range("H126").Value = 3519.2
$ws.Range("I126").Value = 4176.6665
$ws.Range("K126").Value = 12529.9995
$ws.Range("M126").Value = -7589.999500000002
$ws.Range("H129").Value = 902.82355
$ws.Range("J129").Value = 3010.3333
$ws.Range("L129").Value = 9030.999899999999
$ws.Range("N129").Value = -19030.9999
$ws.Range("H131").Value = 1342.5416
$ws.Range("I131").Value = 796.25
$ws.Range("K131").Value = 2388.75
$ws.Range("M131").Value = 2651.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 19777.666
$ws.Range("I57").Value = 15000
$ws.Range("J57").Value = 20374.875
$ws.Range("K57").Value = 15000
$ws.Range("L57").Value = 20374.875
$ws.Range("M57").Value = -14180
$ws.Range("N57").Value = -22014.875
$ws.Range("H122").Value = 2839.7144
$ws.Range("I122").Value = 2356.6667
$ws.Range("J122").Value = 4047.3333
$ws.Range("K122").Value = 7070.000100000001
$ws.Range("L122").Value = 12141.9999
$ws.Range("M122").Value = -4620.000100000001
$ws.Range("N122").Value = -17041.9999
$ws.Range("H126").Value = 2509.3044
$ws.Range("I126").Value = 2488.6667
$ws.Range("J126").Value = 2583.6
$ws.Range("K126").Value = 7466.000100000001
$ws.Range("L126").Value = 7750.799999999999
$ws.Range("M126").Value = -4996.000100000001
$ws.Range("N126").Value = -12690.8
$ws.Range("H129").Value = 92499
$ws.Range("J129").Value = 92499
$ws.Range("L129").Value = 92499
$ws.Range("N129").Value = -102499
$ws.Range("H132").Value = 83361416
$ws.Range("I132").Value = 100018400
$ws.Range("J132").Value = 76507
$ws.Range("K132").Value = 300055200
$ws.Range("L132").Value = 229521
$ws.Range("M132").Value = -300052670
$ws.Range("N132").Value = -234581

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1400.125
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 1533.5
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 1533.5
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -1909.5
$ws.Range("H55").Value = 283.53845
$ws.Range("I55").Value = 217.8
$ws.Range("J55").Value = 324.625
$ws.Range("K55").Value = 217.8
$ws.Range("L55").Value = 324.625
$ws.Range("M55").Value = -44.80000000000001
$ws.Range("N55").Value = -670.625
$ws.Range("H132").Value = 1960.175
$ws.Range("I132").Value = 1817.4
$ws.Range("J132").Value = 2959.6
$ws.Range("K132").Value = 5452.200000000001
$ws.Range("L132").Value = 8878.799999999999
$ws.Range("M132").Value = -2922.200000000001
$ws.Range("N132").Value = -13938.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1779.8
$ws.Range("I96").Value = 966.6667
$ws.Range("K96").Value = 966.6667
$ws.Range("M96").Value = 406.3333
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 2224.6875
$ws.Range("I126").Value = 2092.4614
$ws.Range("J126").Value = 2797.6667
$ws.Range("K126").Value = 6277.3842
$ws.Range("L126").Value = 8393.000100000001
$ws.Range("M126").Value = -3807.3842
$ws.Range("N126").Value = -13333.0001
$ws.Range("H136").Value = 3116.913
$ws.Range("I136").Value = 2659.9302
$ws.Range("J136").Value = 9667
$ws.Range("K136").Value = 7979.790599999999
$ws.Range("L136").Value = 29001
$ws.Range("M136").Value = -5429.790599999999
$ws.Range("N136").Value = -34101
